$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.528881549835205
$ws.Range("B1").Value = 6.952014923095703
$ws.Range("C1").Value = 5.399936199188232
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 3.801798820495605
